$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvestorKyc")

# "Send Kyc Form To User *" column (P) flipped from "No" to "Yes"
# for the first two investor rows.
$ws.Range("P2").Value = "Yes"
$ws.Range("P3").Value = "Yes"

# Active selection moved to P4.
$ws.Range("P4").Select()
